$d = $word.ActiveDocument

# Update the date line (wdReplaceOne = 1, scoped to the whole doc content; only one match exists)
$d.Content.Find.Execute("2023-08-12 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-13 Sunday", 1) | Out-Null

# Update each table cell value in document order (row-major, 20 rows x 5 cols).
# Use wdReplaceOne (1) instead of wdReplaceAll (2) so the replace stays scoped to
# the individual cell range instead of rewriting every matching cell in the table
# (several expressions like "60+11=71" and "85-28=57" repeat with different targets).
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Find.Execute("13+0=13", $true, $false, $false, $false, $false, $true, 1, $false, "0+26=26", 1) | Out-Null
$cell = $t.Cell(1, 2)
$cell.Range.Find.Execute("16-11=5", $true, $false, $false, $false, $false, $true, 1, $false, "93-65=28", 1) | Out-Null
$cell = $t.Cell(1, 3)
$cell.Range.Find.Execute("93-4=89", $true, $false, $false, $false, $false, $true, 1, $false, "60-12=48", 1) | Out-Null
$cell = $t.Cell(1, 4)
$cell.Range.Find.Execute("26+9=35", $true, $false, $false, $false, $false, $true, 1, $false, "92-43=49", 1) | Out-Null
$cell = $t.Cell(1, 5)
$cell.Range.Find.Execute("93-71=22", $true, $false, $false, $false, $false, $true, 1, $false, "74-22=52", 1) | Out-Null

$cell = $t.Cell(2, 1)
$cell.Range.Find.Execute("82+2=84", $true, $false, $false, $false, $false, $true, 1, $false, "74+12=86", 1) | Out-Null
$cell = $t.Cell(2, 2)
$cell.Range.Find.Execute("58+3=61", $true, $false, $false, $false, $false, $true, 1, $false, "21+75=96", 1) | Out-Null
$cell = $t.Cell(2, 3)
$cell.Range.Find.Execute("50-14=36", $true, $false, $false, $false, $false, $true, 1, $false, "0+40=40", 1) | Out-Null
$cell = $t.Cell(2, 4)
$cell.Range.Find.Execute("66+17=83", $true, $false, $false, $false, $false, $true, 1, $false, "11+31=42", 1) | Out-Null
$cell = $t.Cell(2, 5)
$cell.Range.Find.Execute("13+18=31", $true, $false, $false, $false, $false, $true, 1, $false, "76-46=30", 1) | Out-Null

$cell = $t.Cell(3, 1)
$cell.Range.Find.Execute("9+40=49", $true, $false, $false, $false, $false, $true, 1, $false, "68+31=99", 1) | Out-Null
$cell = $t.Cell(3, 2)
$cell.Range.Find.Execute("75+8=83", $true, $false, $false, $false, $false, $true, 1, $false, "66+8=74", 1) | Out-Null
$cell = $t.Cell(3, 3)
$cell.Range.Find.Execute("83+1=84", $true, $false, $false, $false, $false, $true, 1, $false, "18+20=38", 1) | Out-Null
$cell = $t.Cell(3, 4)
$cell.Range.Find.Execute("38-29=9", $true, $false, $false, $false, $false, $true, 1, $false, "14+81=95", 1) | Out-Null
$cell = $t.Cell(3, 5)
$cell.Range.Find.Execute("56-15=41", $true, $false, $false, $false, $false, $true, 1, $false, "32+24=56", 1) | Out-Null

$cell = $t.Cell(4, 1)
$cell.Range.Find.Execute("0+90=90", $true, $false, $false, $false, $false, $true, 1, $false, "52+18=70", 1) | Out-Null
$cell = $t.Cell(4, 2)
$cell.Range.Find.Execute("11+66=77", $true, $false, $false, $false, $false, $true, 1, $false, "37+37=74", 1) | Out-Null
$cell = $t.Cell(4, 3)
$cell.Range.Find.Execute("87-19=68", $true, $false, $false, $false, $false, $true, 1, $false, "52-49=3", 1) | Out-Null
$cell = $t.Cell(4, 4)
$cell.Range.Find.Execute("80-17=63", $true, $false, $false, $false, $false, $true, 1, $false, "70+25=95", 1) | Out-Null
$cell = $t.Cell(4, 5)
$cell.Range.Find.Execute("20+34=54", $true, $false, $false, $false, $false, $true, 1, $false, "30+47=77", 1) | Out-Null

$cell = $t.Cell(5, 1)
$cell.Range.Find.Execute("77+12=89", $true, $false, $false, $false, $false, $true, 1, $false, "47-29=18", 1) | Out-Null
$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("57-37=20", $true, $false, $false, $false, $false, $true, 1, $false, "9+75=84", 1) | Out-Null
$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("48+48=96", $true, $false, $false, $false, $false, $true, 1, $false, "53+22=75", 1) | Out-Null
$cell = $t.Cell(5, 4)
$cell.Range.Find.Execute("12+79=91", $true, $false, $false, $false, $false, $true, 1, $false, "81-34=47", 1) | Out-Null
$cell = $t.Cell(5, 5)
$cell.Range.Find.Execute("67-42=25", $true, $false, $false, $false, $false, $true, 1, $false, "16+12=28", 1) | Out-Null

$cell = $t.Cell(6, 1)
$cell.Range.Find.Execute("72-20=52", $true, $false, $false, $false, $false, $true, 1, $false, "67-22=45", 1) | Out-Null
$cell = $t.Cell(6, 2)
$cell.Range.Find.Execute("27-22=5", $true, $false, $false, $false, $false, $true, 1, $false, "19+8=27", 1) | Out-Null
$cell = $t.Cell(6, 3)
$cell.Range.Find.Execute("88-40=48", $true, $false, $false, $false, $false, $true, 1, $false, "1+53=54", 1) | Out-Null
$cell = $t.Cell(6, 4)
$cell.Range.Find.Execute("8+70=78", $true, $false, $false, $false, $false, $true, 1, $false, "23+18=41", 1) | Out-Null
$cell = $t.Cell(6, 5)
$cell.Range.Find.Execute("46-24=22", $true, $false, $false, $false, $false, $true, 1, $false, "30+68=98", 1) | Out-Null

$cell = $t.Cell(7, 1)
$cell.Range.Find.Execute("96-67=29", $true, $false, $false, $false, $false, $true, 1, $false, "31+64=95", 1) | Out-Null
$cell = $t.Cell(7, 2)
$cell.Range.Find.Execute("91-46=45", $true, $false, $false, $false, $false, $true, 1, $false, "88-29=59", 1) | Out-Null
$cell = $t.Cell(7, 3)
$cell.Range.Find.Execute("41+58=99", $true, $false, $false, $false, $false, $true, 1, $false, "68-59=9", 1) | Out-Null
$cell = $t.Cell(7, 4)
$cell.Range.Find.Execute("2+22=24", $true, $false, $false, $false, $false, $true, 1, $false, "39+11=50", 1) | Out-Null
$cell = $t.Cell(7, 5)
$cell.Range.Find.Execute("89+3=92", $true, $false, $false, $false, $false, $true, 1, $false, "9+25=34", 1) | Out-Null

$cell = $t.Cell(8, 1)
$cell.Range.Find.Execute("0+99=99", $true, $false, $false, $false, $false, $true, 1, $false, "99-58=41", 1) | Out-Null
$cell = $t.Cell(8, 2)
$cell.Range.Find.Execute("67-1=66", $true, $false, $false, $false, $false, $true, 1, $false, "99-39=60", 1) | Out-Null
$cell = $t.Cell(8, 3)
$cell.Range.Find.Execute("26+7=33", $true, $false, $false, $false, $false, $true, 1, $false, "2+66=68", 1) | Out-Null
$cell = $t.Cell(8, 4)
$cell.Range.Find.Execute("96+2=98", $true, $false, $false, $false, $false, $true, 1, $false, "93+3=96", 1) | Out-Null
$cell = $t.Cell(8, 5)
$cell.Range.Find.Execute("16-7=9", $true, $false, $false, $false, $false, $true, 1, $false, "65-27=38", 1) | Out-Null

$cell = $t.Cell(9, 1)
$cell.Range.Find.Execute("67-40=27", $true, $false, $false, $false, $false, $true, 1, $false, "53-37=16", 1) | Out-Null
$cell = $t.Cell(9, 2)
$cell.Range.Find.Execute("69-39=30", $true, $false, $false, $false, $false, $true, 1, $false, "83-5=78", 1) | Out-Null
$cell = $t.Cell(9, 3)
$cell.Range.Find.Execute("89-60=29", $true, $false, $false, $false, $false, $true, 1, $false, "11-10=1", 1) | Out-Null
$cell = $t.Cell(9, 4)
$cell.Range.Find.Execute("64-57=7", $true, $false, $false, $false, $false, $true, 1, $false, "62-1=61", 1) | Out-Null
$cell = $t.Cell(9, 5)
$cell.Range.Find.Execute("3+70=73", $true, $false, $false, $false, $false, $true, 1, $false, "70+12=82", 1) | Out-Null

$cell = $t.Cell(10, 1)
$cell.Range.Find.Execute("46-15=31", $true, $false, $false, $false, $false, $true, 1, $false, "69-68=1", 1) | Out-Null
$cell = $t.Cell(10, 2)
$cell.Range.Find.Execute("17+38=55", $true, $false, $false, $false, $false, $true, 1, $false, "74-16=58", 1) | Out-Null
$cell = $t.Cell(10, 3)
$cell.Range.Find.Execute("51-46=5", $true, $false, $false, $false, $false, $true, 1, $false, "31+47=78", 1) | Out-Null
$cell = $t.Cell(10, 4)
$cell.Range.Find.Execute("98-50=48", $true, $false, $false, $false, $false, $true, 1, $false, "51+28=79", 1) | Out-Null
$cell = $t.Cell(10, 5)
$cell.Range.Find.Execute("77-2=75", $true, $false, $false, $false, $false, $true, 1, $false, "16-6=10", 1) | Out-Null

$cell = $t.Cell(11, 1)
$cell.Range.Find.Execute("75-57=18", $true, $false, $false, $false, $false, $true, 1, $false, "95-54=41", 1) | Out-Null
$cell = $t.Cell(11, 2)
$cell.Range.Find.Execute("31+13=44", $true, $false, $false, $false, $false, $true, 1, $false, "64-0=64", 1) | Out-Null
$cell = $t.Cell(11, 3)
$cell.Range.Find.Execute("44-22=22", $true, $false, $false, $false, $false, $true, 1, $false, "1+85=86", 1) | Out-Null
$cell = $t.Cell(11, 4)
$cell.Range.Find.Execute("19+9=28", $true, $false, $false, $false, $false, $true, 1, $false, "15+66=81", 1) | Out-Null
$cell = $t.Cell(11, 5)
$cell.Range.Find.Execute("8+12=20", $true, $false, $false, $false, $false, $true, 1, $false, "74+5=79", 1) | Out-Null

$cell = $t.Cell(12, 1)
$cell.Range.Find.Execute("95-35=60", $true, $false, $false, $false, $false, $true, 1, $false, "68+26=94", 1) | Out-Null
$cell = $t.Cell(12, 2)
$cell.Range.Find.Execute("7+17=24", $true, $false, $false, $false, $false, $true, 1, $false, "91-71=20", 1) | Out-Null
$cell = $t.Cell(12, 3)
$cell.Range.Find.Execute("22+36=58", $true, $false, $false, $false, $false, $true, 1, $false, "54-50=4", 1) | Out-Null
$cell = $t.Cell(12, 4)
$cell.Range.Find.Execute("5+56=61", $true, $false, $false, $false, $false, $true, 1, $false, "69-8=61", 1) | Out-Null
$cell = $t.Cell(12, 5)
$cell.Range.Find.Execute("35+47=82", $true, $false, $false, $false, $false, $true, 1, $false, "95-47=48", 1) | Out-Null

$cell = $t.Cell(13, 1)
$cell.Range.Find.Execute("86-35=51", $true, $false, $false, $false, $false, $true, 1, $false, "62-43=19", 1) | Out-Null
$cell = $t.Cell(13, 2)
$cell.Range.Find.Execute("85-28=57", $true, $false, $false, $false, $false, $true, 1, $false, "43-34=9", 1) | Out-Null
$cell = $t.Cell(13, 3)
$cell.Range.Find.Execute("73-39=34", $true, $false, $false, $false, $false, $true, 1, $false, "56+2=58", 1) | Out-Null
$cell = $t.Cell(13, 4)
$cell.Range.Find.Execute("93-80=13", $true, $false, $false, $false, $false, $true, 1, $false, "51-35=16", 1) | Out-Null
$cell = $t.Cell(13, 5)
$cell.Range.Find.Execute("54-24=30", $true, $false, $false, $false, $false, $true, 1, $false, "98-66=32", 1) | Out-Null

$cell = $t.Cell(14, 1)
$cell.Range.Find.Execute("60+11=71", $true, $false, $false, $false, $false, $true, 1, $false, "18+62=80", 1) | Out-Null
$cell = $t.Cell(14, 2)
$cell.Range.Find.Execute("44+24=68", $true, $false, $false, $false, $false, $true, 1, $false, "38+30=68", 1) | Out-Null
$cell = $t.Cell(14, 3)
$cell.Range.Find.Execute("8+66=74", $true, $false, $false, $false, $false, $true, 1, $false, "81-23=58", 1) | Out-Null
$cell = $t.Cell(14, 4)
$cell.Range.Find.Execute("60+11=71", $true, $false, $false, $false, $false, $true, 1, $false, "77-15=62", 1) | Out-Null
$cell = $t.Cell(14, 5)
$cell.Range.Find.Execute("92-33=59", $true, $false, $false, $false, $false, $true, 1, $false, "95-85=10", 1) | Out-Null

$cell = $t.Cell(15, 1)
$cell.Range.Find.Execute("14-2=12", $true, $false, $false, $false, $false, $true, 1, $false, "35+21=56", 1) | Out-Null
$cell = $t.Cell(15, 2)
$cell.Range.Find.Execute("76-64=12", $true, $false, $false, $false, $false, $true, 1, $false, "28+28=56", 1) | Out-Null
$cell = $t.Cell(15, 3)
$cell.Range.Find.Execute("16+33=49", $true, $false, $false, $false, $false, $true, 1, $false, "50+32=82", 1) | Out-Null
$cell = $t.Cell(15, 4)
$cell.Range.Find.Execute("26+72=98", $true, $false, $false, $false, $false, $true, 1, $false, "4+36=40", 1) | Out-Null
$cell = $t.Cell(15, 5)
$cell.Range.Find.Execute("58+15=73", $true, $false, $false, $false, $false, $true, 1, $false, "8+63=71", 1) | Out-Null

$cell = $t.Cell(16, 1)
$cell.Range.Find.Execute("53+42=95", $true, $false, $false, $false, $false, $true, 1, $false, "22-12=10", 1) | Out-Null
$cell = $t.Cell(16, 2)
$cell.Range.Find.Execute("82-79=3", $true, $false, $false, $false, $false, $true, 1, $false, "18+47=65", 1) | Out-Null
$cell = $t.Cell(16, 3)
$cell.Range.Find.Execute("30-21=9", $true, $false, $false, $false, $false, $true, 1, $false, "32-2=30", 1) | Out-Null
$cell = $t.Cell(16, 4)
$cell.Range.Find.Execute("86-58=28", $true, $false, $false, $false, $false, $true, 1, $false, "21+48=69", 1) | Out-Null
$cell = $t.Cell(16, 5)
$cell.Range.Find.Execute("70-20=50", $true, $false, $false, $false, $false, $true, 1, $false, "41+19=60", 1) | Out-Null

$cell = $t.Cell(17, 1)
$cell.Range.Find.Execute("24+33=57", $true, $false, $false, $false, $false, $true, 1, $false, "88-37=51", 1) | Out-Null
$cell = $t.Cell(17, 2)
$cell.Range.Find.Execute("94+1=95", $true, $false, $false, $false, $false, $true, 1, $false, "58-31=27", 1) | Out-Null
$cell = $t.Cell(17, 3)
$cell.Range.Find.Execute("3+41=44", $true, $false, $false, $false, $false, $true, 1, $false, "70+10=80", 1) | Out-Null
$cell = $t.Cell(17, 4)
$cell.Range.Find.Execute("9+46=55", $true, $false, $false, $false, $false, $true, 1, $false, "57+31=88", 1) | Out-Null
$cell = $t.Cell(17, 5)
$cell.Range.Find.Execute("71+6=77", $true, $false, $false, $false, $false, $true, 1, $false, "79-8=71", 1) | Out-Null

$cell = $t.Cell(18, 1)
$cell.Range.Find.Execute("59-39=20", $true, $false, $false, $false, $false, $true, 1, $false, "55+38=93", 1) | Out-Null
$cell = $t.Cell(18, 2)
$cell.Range.Find.Execute("71+12=83", $true, $false, $false, $false, $false, $true, 1, $false, "93-25=68", 1) | Out-Null
$cell = $t.Cell(18, 3)
$cell.Range.Find.Execute("76-54=22", $true, $false, $false, $false, $false, $true, 1, $false, "6+61=67", 1) | Out-Null
$cell = $t.Cell(18, 4)
$cell.Range.Find.Execute("47-42=5", $true, $false, $false, $false, $false, $true, 1, $false, "23+8=31", 1) | Out-Null
$cell = $t.Cell(18, 5)
$cell.Range.Find.Execute("38-23=15", $true, $false, $false, $false, $false, $true, 1, $false, "14+64=78", 1) | Out-Null

$cell = $t.Cell(19, 1)
$cell.Range.Find.Execute("50+8=58", $true, $false, $false, $false, $false, $true, 1, $false, "37+6=43", 1) | Out-Null
$cell = $t.Cell(19, 2)
$cell.Range.Find.Execute("98-19=79", $true, $false, $false, $false, $false, $true, 1, $false, "4+54=58", 1) | Out-Null
$cell = $t.Cell(19, 3)
$cell.Range.Find.Execute("22+20=42", $true, $false, $false, $false, $false, $true, 1, $false, "9-9=0", 1) | Out-Null
$cell = $t.Cell(19, 4)
$cell.Range.Find.Execute("57+38=95", $true, $false, $false, $false, $false, $true, 1, $false, "7+87=94", 1) | Out-Null
$cell = $t.Cell(19, 5)
$cell.Range.Find.Execute("85-28=57", $true, $false, $false, $false, $false, $true, 1, $false, "97-20=77", 1) | Out-Null

$cell = $t.Cell(20, 1)
$cell.Range.Find.Execute("29+17=46", $true, $false, $false, $false, $false, $true, 1, $false, "13+2=15", 1) | Out-Null
$cell = $t.Cell(20, 2)
$cell.Range.Find.Execute("93-47=46", $true, $false, $false, $false, $false, $true, 1, $false, "88-54=34", 1) | Out-Null
$cell = $t.Cell(20, 3)
$cell.Range.Find.Execute("84-13=71", $true, $false, $false, $false, $false, $true, 1, $false, "8-3=5", 1) | Out-Null
$cell = $t.Cell(20, 4)
$cell.Range.Find.Execute("49-43=6", $true, $false, $false, $false, $false, $true, 1, $false, "39+15=54", 1) | Out-Null
$cell = $t.Cell(20, 5)
$cell.Range.Find.Execute("92-6=86", $true, $false, $false, $false, $false, $true, 1, $false, "85-60=25", 1) | Out-Null

